# Update the "取得日時" (retrieved timestamp) column for all data rows
# on the "ランサーズ" sheet to the new run's timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-19 06:35:10"

$ws.Range("A2:A7").Value = $newTimestamp
